$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.187863986473793
$ws.Range("D4").Value = 0.2254367837685516
$ws.Range("D9").Value = 0.8453879391320683
$ws.Range("D13").Value = 0.1502911891790344
$ws.Range("D16").Value = 0.826601540484689
$ws.Range("D20").Value = 0.9862859289874132
$ws.Range("D23").Value = 0.4602667668607928
$ws.Range("D24").Value = 0.1502911891790344
$ws.Range("D27").Value = 0.09393199323689649
$ws.Range("D30").Value = 1.427766297200826
$ws.Range("D31").Value = 0.1127183918842758
$ws.Range("D32").Value = 0.5166259628029307
$ws.Range("D34").Value = 1.315047905316551
$ws.Range("D35").Value = 0.1033251925605861
$ws.Range("D36").Value = 0.9956791283111027
$ws.Range("D40").Value = 0.8641743377794476
$ws.Range("D41").Value = 1.87863986473793
$ws.Range("D43").Value = 0.0187863986473793
$ws.Range("D46").Value = 0.4320871688897238
$ws.Range("D47").Value = 0.1690775878264137
$ws.Range("D53").Value = 0.1221115912079654
$ws.Range("D55").Value = 0.04696599661844825
$ws.Range("D56").Value = 0.187863986473793
$ws.Range("D58").Value = 0.7608491452188616
$ws.Range("D59").Value = 0.3569415743002066
$ws.Range("D61").Value = 0.2442231824159309
$ws.Range("D62").Value = 0.2630095810633102
$ws.Range("D65").Value = 0.5260191621266204
$ws.Range("D67").Value = 0.1033251925605861
$ws.Range("D71").Value = 0.0187863986473793
$ws.Range("D73").Value = 0.4978395641555514
$ws.Range("D76").Value = 0.06575239526582755
$ws.Range("D77").Value = 0.1221115912079654
$ws.Range("D78").Value = 0.1408979898553447
$ws.Range("D79").Value = 1.493518692466654
$ws.Range("D81").Value = 0.1690775878264137
$ws.Range("D82").Value = 0.1972571857974826
$ws.Range("D85").Value = 0.1221115912079654
$ws.Range("D88").Value = 0.2536163817396205
$ws.Range("D92").Value = 0.1033251925605861
$ws.Range("D94").Value = 2.047717452564344
$ws.Range("D99").Value = 0.1033251925605861
$ws.Range("D100").Value = 0.0375727972947586
$ws.Range("D101").Value = 0.6105579560398272
$ws.Range("D107").Value = 2.911891790343791
$ws.Range("D108").Value = 3.024610182228067
$ws.Range("D110").Value = 0.3851211722712756
$ws.Range("D113").Value = 2.864925793725343
$ws.Range("D117").Value = 2.094683449182792
$ws.Range("D118").Value = 4.724779259815894
$ws.Range("D120").Value = 2.66766860792786
$ws.Range("D127").Value = 0.7138831486004132
$ws.Range("D135").Value = 0.4133007702423445
$ws.Range("D136").Value = 0.1033251925605861
$ws.Range("D137").Value = 0.0563591959421379
$ws.Range("D142").Value = 0.08453879391320684
$ws.Range("D147").Value = 0.3851211722712756
$ws.Range("D149").Value = 0.2911891790343791
$ws.Range("D151").Value = 0.2348299830922412
$ws.Range("D152").Value = 0.1408979898553447
$ws.Range("D157").Value = 0.1033251925605861
$ws.Range("D158").Value = 0.4320871688897238
$ws.Range("D161").Value = 0.04696599661844825
$ws.Range("D162").Value = 1.06143152357693
$ws.Range("D166").Value = 0.9768927296637234
$ws.Range("D170").Value = 0.0563591959421379
$ws.Range("D171").Value = 0.4133007702423445
$ws.Range("D177").Value = 0.06575239526582755
$ws.Range("D184").Value = 0.1408979898553447
$ws.Range("D185").Value = 1.230509111403344
$ws.Range("D186").Value = 0.3005823783580688
$ws.Range("D190").Value = 0.3287619763291377
$ws.Range("D197").Value = 0.1221115912079654
$ws.Range("D202").Value = 0.0375727972947586
$ws.Range("D203").Value = 0.09393199323689649
$ws.Range("D204").Value = 0.4226939695660342
$ws.Range("D206").Value = 0.08453879391320684
$ws.Range("D207").Value = 0.2724027803869998
$ws.Range("D211").Value = 0.3381551756528274
$ws.Range("D215").Value = 0.4320871688897238
$ws.Range("D216").Value = 1.841067067443171
$ws.Range("D217").Value = 0.4226939695660342
$ws.Range("D218").Value = 0.5541987600976892
$ws.Range("D220").Value = 0.3099755776817584
$ws.Range("D221").Value = 2.151042645124929
$ws.Range("D222").Value = 0.2160435844448619
$ws.Range("D223").Value = 0.3005823783580688
$ws.Range("D224").Value = 1.831673868119482
$ws.Range("D225").Value = 0.1502911891790344
$ws.Range("D229").Value = 0.04696599661844825
$ws.Range("D236").Value = 0.4320871688897238
$ws.Range("D238").Value = 0.1127183918842758
$ws.Range("D240").Value = 0.02817959797106895
